# SI1_RNAseq_libraries_used.xlsx — add two new "Oxford population" library rows
# (216, 217) plus a new "Notes" column (D) documenting them, mirroring the
# author's manual entry in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 216 -----------------------------------------------------------
# Write A/B/C in this order so new shared-string entries are minted in the
# same order the original authoring session produced them.
$ws.Range("A216").Value = "SRRXXXXXXX"
$ws.Range("B216").Value = "PRJNAXXXXXX"

# Column B carries the "BioProject" cell style (font: Arial 10, black) used
# throughout the sheet (e.g. B215). Copy/PasteSpecial(Formats) reuses that
# existing style rather than minting a new one via Font.Name/Size/Color.
$ws.Range("B215").Copy()
$ws.Range("B216").PasteSpecial($xlPasteFormats)

$ws.Range("C216").Value = "Sexual"

# --- New column header ---------------------------------------------------
$ws.Range("D1").Value = "Notes"

$ws.Range("D216").Value = "Oxford populaiton - SRA in process of submission"

# --- Row 217 (duplicate of 216) ------------------------------------------
$ws.Range("A217").Value = "SRRXXXXXXX"
$ws.Range("B217").Value = "PRJNAXXXXXX"
$ws.Range("B215").Copy()
$ws.Range("B217").PasteSpecial($xlPasteFormats)
$ws.Range("C217").Value = "Sexual"
$ws.Range("D217").Value = "Oxford populaiton - SRA in process of submission"

# --- View state: leave selection on the last-entered cell, scrolled down -
$ws.Range("C216").Select()
$excel.ActiveWindow.ScrollRow = 206
